$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the title text in A1 (style/merge stays as-is)
$ws.Range("A1").Value = ""

# Clear out the dummy data rows (3-8) including values and formatting
$ws.Range("A3:E8").Clear()
